$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-02-29 Thursday" "2024-03-01 Friday"

Replace-Text "43×88=3784" "98×84=8232"
Replace-Text "60×73=4380" "21×97=2037"
Replace-Text "96×41=3936" "91×97=8827"
Replace-Text "18×49=882" "61×27=1647"
Replace-Text "85×70=5950" "27×19=513"

Replace-Text "66×77=5082" "60×87=5220"
Replace-Text "86×62=5332" "20×95=1900"
Replace-Text "72×81=5832" "93×13=1209"
Replace-Text "27×22=594" "67×22=1474"
Replace-Text "74×84=6216" "19×70=1330"

Replace-Text "65×66=4290" "60×37=2220"
Replace-Text "18×45=810" "26×27=702"
Replace-Text "42×45=1890" "28×39=1092"
Replace-Text "54×48=2592" "64×29=1856"
Replace-Text "63×92=5796" "55×79=4345"

Replace-Text "86×73=6278" "18×19=342"
Replace-Text "47×37=1739" "17×68=1156"
Replace-Text "46×35=1610" "14×22=308"
Replace-Text "85×13=1105" "94×96=9024"
Replace-Text "66×53=3498" "43×40=1720"

Replace-Text "25×15=375" "12×31=372"
Replace-Text "74×65=4810" "67×47=3149"
Replace-Text "36×81=2916" "67×52=3484"
Replace-Text "28×75=2100" "41×77=3157"
Replace-Text "79×50=3950" "12×51=612"
